$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 label (Capito -> wheap cream)
$ws.Range("B2").Value = "wheap cream"

# Add new rows of data
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = "laptiono"
$ws.Range("C3").Value = 200

$ws.Range("A4").Value = 4
$ws.Range("B4").Value = "lattee"
$ws.Range("C4").Value = 343

$ws.Range("A5").Value = 5
$ws.Range("B5").Value = "Lazy shit"
$ws.Range("C5").Value = 30

# Match the final selection shown in the saved workbook
$ws.Range("C5").Select() | Out-Null
